$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Ativacao date change 01/01/2012 -> 01/01/2021
$ws.Range("B8").Value = '01/01/2021'
$ws.Range("C8").Value = '01/01/2021'

# Row 10: Objetivos (PT) updated text
$ws.Range("B10").Value = 'Apresentar os fundamentos da Contabilidade Geral, de modo que os alunos possam compreender os conceitos e princípios básicos de Contabilidade. Capacitar na produção e análise de documentos e de registros de transações contábeis. Compreender os sistemas de custeio e sua aplicação. Capacitar a projetar e implantar sistemas de custeio. Capacitar a analisar custos em relação a estratégia de mercado e de produção de uma organização.'
$ws.Range("C10").Value = 'Apresentar os fundamentos da Contabilidade Geral, de modo que os alunos possam compreender os conceitos e princípios básicos de Contabilidade. Capacitar na produção e análise de documentos e de registros de transações contábeis. Compreender os sistemas de custeio e sua aplicação. Capacitar a projetar e implantar sistemas de custeio. Capacitar a analisar custos em relação a estratégia de mercado e de produção de uma organização.'

# Row 11: Objectives (EN) - new cells, copy format from row 10
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = 'Present the fundamentals of General Accounting, so that students can understand the basic concepts and principles of Accounting. Train in the production and analysis of documents and records of accounting transactions. Understand costing systems and their application. Empower to design and implement costing systems. Empower to analyze costs in relation to an organization''s market and production strategy.'
$ws.Range("C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 'Present the fundamentals of General Accounting, so that students can understand the basic concepts and principles of Accounting. Train in the production and analysis of documents and records of accounting transactions. Understand costing systems and their application. Empower to design and implement costing systems. Empower to analyze costs in relation to an organization''s market and production strategy.'

# Row 13: Docente responsavel updated
$ws.Range("B13").Value = '11079086 - Herlandí de Souza Andrade'
$ws.Range("C13").Value = '11079086 - Herlandí de Souza Andrade'

# Row 15: Short syllabus (EN) - new cells, copy format from row 14
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = '1 - Fundamentals of Accounting. 2 - Costing Systems. 3 - Costing Systems Project. 4 - Cost Analysis.'
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = '1 - Fundamentals of Accounting. 2 - Costing Systems. 3 - Costing Systems Project. 4 - Cost Analysis.'

# Row 16: Programa (PT) - text reflowed (no internal newlines)
$ws.Range("B16").Value = '1. FUNDAMENTOS DE CONTABILIDADEConceito de contabilidade. Contabilidade e Engenharia de Produção. Balanço: ativo, passivo, patrimônio líquido. Procedimentos contábeis básicos. Variações da situação líquida; despesa, receita. Regimes de competência e caixa. Receitas e despesas diferidas. Fatos Contábeis. Operações com mercadorias. Inventário. Ativo Imobilizado e Amortização. Demonstrativo de origens e aplicações. Análise de Balanço e de resultados.2. SISTEMAS DE CUSTEIOCustos diretos e indiretos, fixos e variáveis. Acumulação de custos, classificação, fatores de custo. Custeio por Absorção. Produção por ordem, contínua, conjunta. Custeio Direto: margem de contribuição. Custeio ABC.3. PROJETO DE SISTEMAS DE CUSTEIOProdutos e Departamentos. Sistemas de produção e sistemas de custeio. O problema da inflação.4. ANÁLISE DE CUSTOSCusto-Volume-Lucro. Contribuição marginal. Análise de variações. Equação de produtividade global. Alavancagem operacional. TIR e lucratividade.'
$ws.Range("C16").Value = '1. FUNDAMENTOS DE CONTABILIDADEConceito de contabilidade. Contabilidade e Engenharia de Produção. Balanço: ativo, passivo, patrimônio líquido. Procedimentos contábeis básicos. Variações da situação líquida; despesa, receita. Regimes de competência e caixa. Receitas e despesas diferidas. Fatos Contábeis. Operações com mercadorias. Inventário. Ativo Imobilizado e Amortização. Demonstrativo de origens e aplicações. Análise de Balanço e de resultados.2. SISTEMAS DE CUSTEIOCustos diretos e indiretos, fixos e variáveis. Acumulação de custos, classificação, fatores de custo. Custeio por Absorção. Produção por ordem, contínua, conjunta. Custeio Direto: margem de contribuição. Custeio ABC.3. PROJETO DE SISTEMAS DE CUSTEIOProdutos e Departamentos. Sistemas de produção e sistemas de custeio. O problema da inflação.4. ANÁLISE DE CUSTOSCusto-Volume-Lucro. Contribuição marginal. Análise de variações. Equação de produtividade global. Alavancagem operacional. TIR e lucratividade.'

# Row 17: Syllabus (EN) - new cells, copy format from row 16
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = '1. ACCOUNTING FUNDAMENTALSAccounting concept. Production Accounting and Engineering. Balance sheet: assets, liabilities, equity. Basic accounting procedures. Changes in equity; expense, revenue. Accrual and cash regimes. Deferred income and expenses. Accounting Facts. Goods operations. Inventory. Property, Plant and Equipment and Amortization. Statement of origins and applications. Balance sheet and results analysis.2. COSTING SYSTEMSDirect and indirect costs, fixed and variable. Accumulation of costs, classification, cost factors. Absorption costing. Orderly, continuous, joint production. Direct Costing: contribution margin. ABC costing.3. COSTING SYSTEMS DESIGNProducts and Departments. Production and costing systems. The inflation problem.4. COST ANALYSISCost-Volume-Profit. Marginal contribution. Analysis of variations. Global productivity equation. Operational leverage. IRR and profitability.'
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = '1. ACCOUNTING FUNDAMENTALSAccounting concept. Production Accounting and Engineering. Balance sheet: assets, liabilities, equity. Basic accounting procedures. Changes in equity; expense, revenue. Accrual and cash regimes. Deferred income and expenses. Accounting Facts. Goods operations. Inventory. Property, Plant and Equipment and Amortization. Statement of origins and applications. Balance sheet and results analysis.2. COSTING SYSTEMSDirect and indirect costs, fixed and variable. Accumulation of costs, classification, cost factors. Absorption costing. Orderly, continuous, joint production. Direct Costing: contribution margin. ABC costing.3. COSTING SYSTEMS DESIGNProducts and Departments. Production and costing systems. The inflation problem.4. COST ANALYSISCost-Volume-Profit. Marginal contribution. Analysis of variations. Global productivity equation. Operational leverage. IRR and profitability.'

# Row 19: Metodo updated
$ws.Range("B19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'
$ws.Range("C19").Value = 'Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras.'

# Row 20: Criterio updated
$ws.Range("B20").Value = 'Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'
$ws.Range("C20").Value = 'Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.'

# Row 21: Norma de recuperacao updated
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.'

# Row 22: Bibliografia updated (reflowed, no internal newlines)
$ws.Range("B22").Value = 'IUDICIBUS, S.; MARION, J. C. Curso de Contabilidade para não Contadores. 8 ed. São Paulo: Atlas, 2018.MARTINS, E. Contabilidade de Custos. 11 ed. São Paulo: Atlas, 2018.SANTOS, J. J. Manual de Contabilidade e Análise de Custos. 7 ed. São Paulo: Atlas, 2017.SILVA, R. N. S.; LINS, L. S. Gestão de Custos - Contabilidade, Controle e Análise. 4 ed. São Paulo: Atlas, 2017.DUTRA, R. G. Custos: Uma Abordagem Prática. 8 ed. São Paulo: Atlas, 2017.MEGLIORINI, E.; BUENO, A. S. Contabilidade para cursos de Engenharia. São Paulo: Atlas, 2014.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.Bruni, Adriano L., Fama, Rubens. Gestão de custos e formação de preços. 5. Ed., São Paulo: Atlas, 2008.Hansen, Don R., Mowen, Maryanne M., Gestão de custos – contabilidade e controle. São Paulo: Pioneira Thomson, 2001Horngren, Charles T., Datar, Srikant M., Foster, George. Contabilidade de custos (vol. 1 e 2). 11. Ed. São Paulo: Pearson, 2004.Maher, Michael. Contabilidade de custos – criando valor para a administração. 5. ed., São Paulo: Atlas, 2001.Martins, Eliseu. Contabilidade de Custos (livro de exercícios). 9. Ed., São Paulo: Atlas, 2006.Martins, Eliseu. Contabilidade de Custos (livro texto). 9. Ed., São Paulo: Atlas, 2003.'
$ws.Range("C22").Value = 'IUDICIBUS, S.; MARION, J. C. Curso de Contabilidade para não Contadores. 8 ed. São Paulo: Atlas, 2018.MARTINS, E. Contabilidade de Custos. 11 ed. São Paulo: Atlas, 2018.SANTOS, J. J. Manual de Contabilidade e Análise de Custos. 7 ed. São Paulo: Atlas, 2017.SILVA, R. N. S.; LINS, L. S. Gestão de Custos - Contabilidade, Controle e Análise. 4 ed. São Paulo: Atlas, 2017.DUTRA, R. G. Custos: Uma Abordagem Prática. 8 ed. São Paulo: Atlas, 2017.MEGLIORINI, E.; BUENO, A. S. Contabilidade para cursos de Engenharia. São Paulo: Atlas, 2014.MORANTE, A. S. Análise das Demonstrações Financeiras. 2 ed. São Paulo: Atlas, 2009.Bruni, Adriano L., Fama, Rubens. Gestão de custos e formação de preços. 5. Ed., São Paulo: Atlas, 2008.Hansen, Don R., Mowen, Maryanne M., Gestão de custos – contabilidade e controle. São Paulo: Pioneira Thomson, 2001Horngren, Charles T., Datar, Srikant M., Foster, George. Contabilidade de custos (vol. 1 e 2). 11. Ed. São Paulo: Pearson, 2004.Maher, Michael. Contabilidade de custos – criando valor para a administração. 5. ed., São Paulo: Atlas, 2001.Martins, Eliseu. Contabilidade de Custos (livro de exercícios). 9. Ed., São Paulo: Atlas, 2006.Martins, Eliseu. Contabilidade de Custos (livro texto). 9. Ed., São Paulo: Atlas, 2003.'

Write-Host "Edit complete"
